$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") stores a date serial number that the automatic
# update bumped by one day (2023-09-08 -> 2023-09-09, i.e. 45177 -> 45178)
# for every data row. Determine the last used data row dynamically and then
# set the whole C2:C<lastRow> range in one shot (Excel broadcasts the scalar
# value to every cell in the range).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

if ($lastRow -ge 2) {
    $rng = $ws.Range("C2:C$lastRow")
    $rng.Value2 = 45178
}
